# Excel file w/leaning scores for news sites
#
# Adds a "Score (v2)" leaning-score column (C) entry for several news
# sites, rescales some existing scores from +/-1 to +/-0.5, replaces the
# "msn" entry in row 22 with "washingtonpost", and appends two new news
# sites (pbs, economist) with their scores in rows 32-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rescale existing scores from +/-1 to +/-0.5 ---
$ws.Range("C3").Value  = -0.5   # ap
$ws.Range("C4").Value  = 0.5    # wsj
$ws.Range("C5").Value  = -0.5   # ft
$ws.Range("C8").Value  = -0.5   # abc
$ws.Range("C9").Value  = -0.5   # cbs
$ws.Range("C18").Value = -0.5   # bloomberg
$ws.Range("C19").Value = -0.5   # reuters
$ws.Range("C21").Value = -0.5   # forbes
$ws.Range("C26").Value = -0.5   # bbc
$ws.Range("C28").Value = -0.5   # cnbc

# --- Fill in previously-empty scores ---
$ws.Range("C11").Value = -1     # nytimes
$ws.Range("C12").Value = -0.5   # usatoday
$ws.Range("C13").Value = -1     # vice
$ws.Range("C15").Value = -0.5   # npr
$ws.Range("C17").Value = -0.5   # businessinsider
$ws.Range("C22").Value = -1     # washingtonpost (see rename below)
$ws.Range("C24").Value = -0.5   # thehill
$ws.Range("C30").Value = -0.5   # news.sky

# --- Row 22: rename news site from "msn" to "washingtonpost" ---
$ws.Range("B22").Value = "washingtonpost"

# --- New rows: additional news sites with their scores ---
$ws.Range("B32").Value = "pbs"
$ws.Range("C32").Value = -0.5
$ws.Range("B33").Value = "economist"
$ws.Range("C33").Value = -0.5

# --- Update the saved view/selection state ---
$ws.Range("C39").Select()
